$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 24.35712233333334
$ws.Range("H2").Value = 73.07136700000001
$ws.Range("I2").Value = 0.3750500562097488
$ws.Range("J2").Value = 0.3750500562097488
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.103724333333333
$ws.Range("N2").Value = 3.311173
$ws.Range("O2").Value = 0.01870879385910814
$ws.Range("P2").Value = 0.01870879385910814
$ws.Range("Q2").Value = 26.88354860927678
$ws.Range("R2").Value = 241.951937483491
$ws.Range("S2").Value = 0.007016734188475113
$ws.Range("T2").Value = 0.007016734188475112

# Row 3
$ws.Range("G3").Value = 24.35712233333334
$ws.Range("H3").Value = 73.07136700000001
$ws.Range("I3").Value = 0.3750500562097488
$ws.Range("J3").Value = 0.3750500562097488
$ws.Range("O3").Value = 0.1603368629650925
$ws.Range("P3").Value = 0.1603368629650925
$ws.Range("Q3").Value = 230.3956033639517
$ws.Range("R3").Value = 2073.560430275565
$ws.Range("S3").Value = 0.06013434946755275
$ws.Range("T3").Value = 0.06013434946755274

# Row 4
$ws.Range("G4").Value = 24.35712233333334
$ws.Range("H4").Value = 73.07136700000001
$ws.Range("I4").Value = 0.3750500562097488
$ws.Range("J4").Value = 0.3750500562097488
$ws.Range("M4").Value = 47.61312599999999
$ws.Range("N4").Value = 142.839378
$ws.Range("O4").Value = 0.80707123365805
$ws.Range("P4").Value = 0.80707123365805
$ws.Range("Q4").Value = 1159.718734654414
$ws.Range("R4").Value = 10437.46861188973
$ws.Range("S4").Value = 0.302692111548723
$ws.Range("T4").Value = 0.3026921115487229

# Row 5
$ws.Range("G5").Value = 24.35712233333334
$ws.Range("H5").Value = 73.07136700000001
$ws.Range("I5").Value = 0.3750500562097488
$ws.Range("J5").Value = 0.3750500562097488
$ws.Range("M5").Value = 0.8190333333333334
$ws.Range("N5").Value = 2.4571
$ws.Range("O5").Value = 0.01388310951774934
$ws.Range("P5").Value = 0.01388310951774934
$ws.Range("Q5").Value = 19.94929509507778
$ws.Range("R5").Value = 179.5436558557
$ws.Range("S5").Value = 0.005206861004997987
$ws.Range("T5").Value = 0.005206861004997987

# Row 6
$ws.Range("I6").Value = 0.2805618708302703
$ws.Range("J6").Value = 0.2805618708302702
$ws.Range("K6").Value = 3
$ws.Range("L6").Value = 1
$ws.Range("M6").Value = 1.103724333333333
$ws.Range("N6").Value = 3.311173
$ws.Range("O6").Value = 0.01870879385910814
$ws.Range("P6").Value = 0.01870879385910814
$ws.Range("Q6").Value = 20.11064541250733
$ws.Range("R6").Value = 180.995808712566
$ws.Range("S6").Value = 0.005248974206089253
$ws.Range("T6").Value = 0.005248974206089252

# Row 7
$ws.Range("I7").Value = 0.2805618708302703
$ws.Range("J7").Value = 0.2805618708302702
$ws.Range("O7").Value = 0.1603368629650925
$ws.Range("P7").Value = 0.1603368629650925
$ws.Range("S7").Value = 0.04498441023654304
$ws.Range("T7").Value = 0.04498441023654303

# Row 8
$ws.Range("I8").Value = 0.2805618708302703
$ws.Range("J8").Value = 0.2805618708302702
$ws.Range("M8").Value = 47.61312599999999
$ws.Range("N8").Value = 142.839378
$ws.Range("O8").Value = 0.80707123365805
$ws.Range("P8").Value = 0.80707123365805
$ws.Range("Q8").Value = 867.5451514919639
$ws.Range("R8").Value = 7807.906363427675
$ws.Range("S8").Value = 0.2264334152083967
$ws.Range("T8").Value = 0.2264334152083967

# Row 9
$ws.Range("I9").Value = 0.2805618708302703
$ws.Range("J9").Value = 0.2805618708302702
$ws.Range("M9").Value = 0.8190333333333334
$ws.Range("N9").Value = 2.4571
$ws.Range("O9").Value = 0.01388310951774934
$ws.Range("P9").Value = 0.01388310951774934
$ws.Range("Q9").Value = 14.92337212313334
$ws.Range("R9").Value = 134.3103491082
$ws.Range("S9").Value = 0.003895071179241285
$ws.Range("T9").Value = 0.003895071179241285

# Row 10
$ws.Range("G10").Value = 22.31748066666667
$ws.Range("H10").Value = 66.952442
$ws.Range("I10").Value = 0.3436437303202491
$ws.Range("J10").Value = 0.343643730320249
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 1.103724333333333
$ws.Range("N10").Value = 3.311173
$ws.Range("O10").Value = 0.01870879385910814
$ws.Range("P10").Value = 0.01870879385910814
$ws.Range("Q10").Value = 24.63234647049622
$ws.Range("R10").Value = 221.691118234466
$ws.Range("S10").Value = 0.006429159711536491
$ws.Range("T10").Value = 0.00642915971153649

# Row 11
$ws.Range("G11").Value = 22.31748066666667
$ws.Range("H11").Value = 66.952442
$ws.Range("I11").Value = 0.3436437303202491
$ws.Range("J11").Value = 0.343643730320249
$ws.Range("O11").Value = 0.1603368629650925
$ws.Range("P11").Value = 0.1603368629650925
$ws.Range("Q11").Value = 211.1025002622434
$ws.Range("R11").Value = 1899.92250236019
$ws.Range("S11").Value = 0.05509875769717099
$ws.Range("T11").Value = 0.05509875769717098

# Row 12
$ws.Range("G12").Value = 22.31748066666667
$ws.Range("H12").Value = 66.952442
$ws.Range("I12").Value = 0.3436437303202491
$ws.Range("J12").Value = 0.343643730320249
$ws.Range("M12").Value = 47.61312599999999
$ws.Range("N12").Value = 142.839378
$ws.Range("O12").Value = 0.80707123365805
$ws.Range("P12").Value = 0.80707123365805
$ws.Range("Q12").Value = 1062.605018984564
$ws.Range("R12").Value = 9563.445170861076
$ws.Range("S12").Value = 0.2773449693684177
$ws.Range("T12").Value = 0.2773449693684176

# Row 13
$ws.Range("G13").Value = 22.31748066666667
$ws.Range("H13").Value = 66.952442
$ws.Range("I13").Value = 0.3436437303202491
$ws.Range("J13").Value = 0.343643730320249
$ws.Range("M13").Value = 0.8190333333333334
$ws.Range("N13").Value = 2.4571
$ws.Range("O13").Value = 0.01388310951774934
$ws.Range("P13").Value = 0.01388310951774934
$ws.Range("Q13").Value = 18.27876058202223
$ws.Range("R13").Value = 164.5088452382
$ws.Range("S13").Value = 0.004770843543123936
$ws.Range("T13").Value = 0.004770843543123935

# Row 14
$ws.Range("E14").Value = 2
$ws.Range("F14").Value = 0.6666666666666666
$ws.Range("G14").Value = 0.04834033333333334
$ws.Range("H14").Value = 0.145021
$ws.Range("I14").Value = 0.0007443426397318391
$ws.Range("J14").Value = 0.0007443426397318388
$ws.Range("K14").Value = 3
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = 1.103724333333333
$ws.Range("N14").Value = 3.311173
$ws.Range("O14").Value = 0.01870879385910814
$ws.Range("P14").Value = 0.01870879385910814
$ws.Range("Q14").Value = 0.05335440218144445
$ws.Range("R14").Value = 0.4801896196330001
$ws.Range("S14").Value = 0.00001392575300728738
$ws.Range("T14").Value = 0.00001392575300728737

# Row 15
$ws.Range("E15").Value = 2
$ws.Range("F15").Value = 0.6666666666666666
$ws.Range("G15").Value = 0.04834033333333334
$ws.Range("H15").Value = 0.145021
$ws.Range("I15").Value = 0.0007443426397318391
$ws.Range("J15").Value = 0.0007443426397318388
$ws.Range("O15").Value = 0.1603368629650925
$ws.Range("P15").Value = 0.1603368629650925
$ws.Range("Q15").Value = 0.4572543551216667
$ws.Range("R15").Value = 4.115289196095
$ws.Range("S15").Value = 0.0001193455638257591
$ws.Range("T15").Value = 0.0001193455638257591

# Row 16
$ws.Range("E16").Value = 2
$ws.Range("F16").Value = 0.6666666666666666
$ws.Range("G16").Value = 0.04834033333333334
$ws.Range("H16").Value = 0.145021
$ws.Range("I16").Value = 0.0007443426397318391
$ws.Range("J16").Value = 0.0007443426397318388
$ws.Range("M16").Value = 47.61312599999999
$ws.Range("N16").Value = 142.839378
$ws.Range("O16").Value = 0.80707123365805
$ws.Range("P16").Value = 0.80707123365805
$ws.Range("Q16").Value = 2.301634381882
$ws.Range("R16").Value = 20.714709436938
$ws.Range("S16").Value = 0.0006007375325126648
$ws.Range("T16").Value = 0.0006007375325126647

# Row 17
$ws.Range("E17").Value = 2
$ws.Range("F17").Value = 0.6666666666666666
$ws.Range("G17").Value = 0.04834033333333334
$ws.Range("H17").Value = 0.145021
$ws.Range("I17").Value = 0.0007443426397318391
$ws.Range("J17").Value = 0.0007443426397318388
$ws.Range("M17").Value = 0.8190333333333334
$ws.Range("N17").Value = 2.4571
$ws.Range("O17").Value = 0.01388310951774934
$ws.Range("P17").Value = 0.01388310951774934
$ws.Range("Q17").Value = 0.03959234434444445
$ws.Range("R17").Value = 0.3563310991
$ws.Range("S17").Value = 0.00001033379038612776
$ws.Range("T17").Value = 0.00001033379038612776
